$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand new row before row 178, shifting rows 178-246 down to 179-247.
$ws.Rows.Item(178).Insert()

# Populate the newly inserted row 178 with its data.
$ws.Range("A178").Value = 10
$ws.Range("B178").Value = "Vega Modelo de Temuco"
$ws.Range("C178").Value = "La Araucanía"
$ws.Range("D178").Value = 44845
$ws.Range("E178").Value = 9
$ws.Range("F178").Value = 100112005
$ws.Range("G178").Value = "Puerro"
$ws.Range("H178").Value = "Sin especificar"
$ws.Range("I178").Value = "Primera"
$ws.Range("J178").Value = 55
$ws.Range("K178").Value = 10000
$ws.Range("L178").Value = 10000
$ws.Range("M178").Value = 10000
$ws.Range("N178").Value = "$/docena de paquetes"
$ws.Range("O178").Value = "Región del Maule"
$ws.Range("P178").Value = 833
$ws.Range("Q178").Value = 12
$ws.Range("R178").Value = "Hortaliza"
